$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.910.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.602.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "554.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.19%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.14%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.595"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.82"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.81%  "

$ws.Range("E10").Value = "  -1.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.142"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.335"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.11%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.062.91"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "58.880.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.87"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.609.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "

$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("E18").Value = "  -0.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.27%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.17"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.74%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("E24").Value = "  +2.18%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.160"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.67%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0756"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  +1.17%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "153.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.73%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.99"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.76%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.894"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.41%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "

$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.871"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.03"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("E39").Value = "  +0.62%  "

$ws.Range("E40").Value = "  -0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "282.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.08%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.600"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.01%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0537"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0953"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0228"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.944.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.12%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.51"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "117.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.97"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.78%  "

